{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in specific resume bullet paragraphs by bolding them and giving\n// them the corporate-blue accent color (#2C3E50). Each target paragraph is\n// located by an EXACT full-text match (so lookalike/decoy paragraphs that\n// share substrings, e.g. \"Impact: Reduced mapping costs by 73.5%...\", are\n// left untouched), then the metric substrings inside it are located with\n// Range.search() and bolded/colored in place.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Each entry: the paragraph's exact current text, and the list of metric\n// substrings (in left-to-right order) that must become bold + colored.\nconst edits = [\n  {\n    text: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    metrics: [\"23%\", \"64%\"]\n  },\n  {\n    text: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00b14.2% to \\u00b12.1%\",\n    metrics: [\"87%\", \"71%\", \"\\u00b14.2%\", \"\\u00b12.1%\"]\n  },\n  {\n    text: \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    metrics: [\"1,200\"]\n  },\n  {\n    text: \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    metrics: [\"$400M\", \"$1B\"]\n  },\n  {\n    text: \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    metrics: [\"73.5%\", \"$4.7M\"]\n  },\n  {\n    text: \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    metrics: [\"87%\", \"71%\"]\n  }\n];\n\n// Track how many times each exact text has already been matched/consumed so\n// paragraphs with duplicate text (none here, but defensive) are each edited\n// exactly once and in document order.\nconst consumed = new Array(edits.length).fill(0);\n\nfor (const paragraph of paragraphs.items) {\n  const editIndex = edits.findIndex((e, i) => e.text === paragraph.text && consumed[i] === 0);\n  if (editIndex === -1) continue;\n  consumed[editIndex] = 1;\n  const metrics = edits[editIndex].metrics;\n\n  for (const metric of metrics) {\n    const hits = paragraph.search(metric, { matchCase: true });\n    hits.load(\"items\");\n    await context.sync();\n    for (const hit of hits.items) {\n      hit.font.bold = true;\n      hit.font.color = \"2C3E50\";\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in specific resume bullet paragraphs by bolding them and giving\n# them the corporate-blue accent color (#2C3E50). Each target paragraph is\n# located by an EXACT full-text match (so lookalike/decoy paragraphs that\n# share substrings, e.g. \"Impact: Reduced mapping costs by 73.5%...\", are\n# left untouched), then the metric substrings inside it are located with\n# Find.Execute() (scoped to that paragraph's Range) and bolded/colored.\n\n$d = $word.ActiveDocument\n\n$bullet = [char]0x2022\n$pm = [char]0x00B1\n\n# Word's Font.Color is a BGR-packed long (wdColor), not an RGB hex triple \u2014\n# build it from the target hex 2C3E50 (R=2C G=3E B=50) so the intent stays\n# legible instead of a magic constant.\n$accentColor = 0x2C + (0x3E * 256) + (0x50 * 65536)\n\n$targets = @(\n    @{\n        Text = \"$bullet Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Metrics = @(\"23%\", \"64%\")\n    },\n    @{\n        Text = \"$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ${pm}4.2% to ${pm}2.1%\"\n        Metrics = @(\"87%\", \"71%\", \"${pm}4.2%\", \"${pm}2.1%\")\n    },\n    @{\n        Text = \"$bullet Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Metrics = @(\"1,200\")\n    },\n    @{\n        Text = \"$bullet Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Metrics = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Text = \"$bullet Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Metrics = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text = \"$bullet Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Metrics = @(\"87%\", \"71%\")\n    }\n)\n\n# Track consumption per target so each exact-text paragraph is only edited\n# once even if (hypothetically) duplicated, and matches happen in document\n# order.\n$consumed = @{}\n\nforeach ($p in $d.Paragraphs) {\n    $ptext = $p.Range.Text.TrimEnd([char]13)\n\n    for ($ti = 0; $ti -lt $targets.Count; $ti++) {\n        if ($consumed.ContainsKey($ti)) { continue }\n        if ($ptext -ne $targets[$ti].Text) { continue }\n\n        $consumed[$ti] = $true\n\n        foreach ($metric in $targets[$ti].Metrics) {\n            $rng = $p.Range\n            $find = $rng.Find\n            $find.ClearFormatting()\n            $find.Text = $metric\n            $find.MatchCase = $true\n            $find.MatchWildcards = $false\n            $find.Forward = $true\n            $find.Wrap = 0\n            $found = $find.Execute()\n            if ($found) {\n                $rng.Font.Bold = $true\n                $rng.Font.Color = $accentColor\n            }\n        }\n\n        break\n    }\n}\n"}
